$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Cells.Item(17, 8).Value = 1807384.2  # H17: was 1692030.8
$ws.Cells.Item(17, 10).Value = 1849377.8  # J17: was 1728777.9
$ws.Cells.Item(17, 12).Value = 5548133.4  # L17: was 5186333.699999999
$ws.Cells.Item(17, 14).Value = -5548469.4  # N17: was -5186669.699999999

# Row 76
$ws.Cells.Item(76, 8).Value = 4974.6  # H76: was 4970.5

# Row 79
$ws.Cells.Item(79, 8).Value = 4974.6  # H79: was 4970.5

# Row 115
$ws.Cells.Item(115, 8).Value = 3926.7  # H115: was 4278.1
$ws.Cells.Item(115, 9).Value = 2374.111  # I115: was 2764.5557
$ws.Cells.Item(115, 11).Value = 7122.333  # K115: was 8293.667099999999
$ws.Cells.Item(115, 13).Value = -5555.333  # M115: was -6726.667099999999

# Row 137
$ws.Cells.Item(137, 8).Value = 4341.8  # H137: was 2489.625
$ws.Cells.Item(137, 9).Value = 4552.25  # I137: was 2350.0476
$ws.Cells.Item(137, 10).Value = 3500  # J137: was 3466.6667
$ws.Cells.Item(137, 11).Value = 13656.75  # K137: was 7050.1428
$ws.Cells.Item(137, 12).Value = 10500  # L137: was 10400.0001
$ws.Cells.Item(137, 13).Value = -11106.75  # M137: was -4500.1428
$ws.Cells.Item(137, 14).Value = -15600  # N137: was -15500.0001

# Row 138
$ws.Cells.Item(138, 8).Value = 290963.53  # H138: was 298746.47
$ws.Cells.Item(138, 10).Value = 1430862  # J138: was 1668839.9
$ws.Cells.Item(138, 12).Value = 4292586  # L138: was 5006519.699999999
$ws.Cells.Item(138, 14).Value = -4302866  # N138: was -5016799.699999999

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 1202.5652  # H2: was 1214.3684
$ws.Cells.Item(2, 9).Value = 1124.2222  # I2: was 1153
$ws.Cells.Item(2, 10).Value = 1484.6  # J2: was 1541.6666
$ws.Cells.Item(2, 11).Value = 1124.2222  # K2: was 1153
$ws.Cells.Item(2, 12).Value = 1484.6  # L2: was 1541.6666
$ws.Cells.Item(2, 13).Value = -1011.2222  # M2: was -1040
$ws.Cells.Item(2, 14).Value = -1710.6  # N2: was -1767.6666

# Row 32
$ws.Cells.Item(32, 8).Value = 9386  # H32: was 9803.666999999999
$ws.Cells.Item(32, 9).Value = 9469.421  # I32: was 9939.277
$ws.Cells.Item(32, 10).Value = 8989.75  # J32: was 8990
$ws.Cells.Item(32, 11).Value = 9469.421  # K32: was 9939.277
$ws.Cells.Item(32, 12).Value = 8989.75  # L32: was 8990
$ws.Cells.Item(32, 13).Value = -9182.421  # M32: was -9652.277
$ws.Cells.Item(32, 14).Value = -9563.75  # N32: was -9564

# Row 61
$ws.Cells.Item(61, 8).Value = 4731.439  # H61: was 4998.3687
$ws.Cells.Item(61, 9).Value = 2564.1614  # I61: was 2694.2144
$ws.Cells.Item(61, 11).Value = 2564.1614  # K61: was 2694.2144
$ws.Cells.Item(61, 13).Value = -2352.1614  # M61: was -2482.2144

# Row 80
$ws.Cells.Item(80, 8).Value = 60296.89  # H80: was 63403.145
$ws.Cells.Item(80, 9).Value = 40100  # I80: was 0
$ws.Cells.Item(80, 10).Value = 62821.5  # J80: was 63403.145
$ws.Cells.Item(80, 11).Value = 40100  # K80: was 0
$ws.Cells.Item(80, 12).Value = 62821.5  # L80: was 63403.145
$ws.Cells.Item(80, 13).Value = -39102  # M80: was None
$ws.Cells.Item(80, 14).Value = -64817.5  # N80: was -65399.145

# Row 83
$ws.Cells.Item(83, 8).Value = 60296.89  # H83: was 63403.145
$ws.Cells.Item(83, 9).Value = 40100  # I83: was 0
$ws.Cells.Item(83, 10).Value = 62821.5  # J83: was 63403.145
$ws.Cells.Item(83, 11).Value = 120300  # K83: was 0
$ws.Cells.Item(83, 12).Value = 188464.5  # L83: was 190209.435
$ws.Cells.Item(83, 13).Value = -115308  # M83: was None
$ws.Cells.Item(83, 14).Value = -198448.5  # N83: was -200193.435

# Row 116
$ws.Cells.Item(116, 8).Value = 1202.5652  # H116: was 1214.3684
$ws.Cells.Item(116, 9).Value = 1124.2222  # I116: was 1153
$ws.Cells.Item(116, 10).Value = 1484.6  # J116: was 1541.6666
$ws.Cells.Item(116, 11).Value = 1124.2222  # K116: was 1153
$ws.Cells.Item(116, 12).Value = 1484.6  # L116: was 1541.6666
$ws.Cells.Item(116, 13).Value = 1169.7778  # M116: was 1141
$ws.Cells.Item(116, 14).Value = -6072.6  # N116: was -6129.6666

# Row 122
$ws.Cells.Item(122, 8).Value = 1805.9756  # H122: was 1759.5333
$ws.Cells.Item(122, 9).Value = 1637.6842  # I122: was 1603.9524
$ws.Cells.Item(122, 11).Value = 4913.0526  # K122: was 4811.857199999999
$ws.Cells.Item(122, 13).Value = -2463.0526  # M122: was -2361.857199999999

# Row 132
$ws.Cells.Item(132, 8).Value = 2423.652  # H132: was 1827.9736
$ws.Cells.Item(132, 9).Value = 1855.0555  # I132: was 1430.909
$ws.Cells.Item(132, 10).Value = 4470.6  # J132: was 4448.6
$ws.Cells.Item(132, 11).Value = 5565.166499999999  # K132: was 4292.727000000001
$ws.Cells.Item(132, 12).Value = 13411.8  # L132: was 13345.8
$ws.Cells.Item(132, 13).Value = -3035.166499999999  # M132: was -1762.727000000001
$ws.Cells.Item(132, 14).Value = -18471.8  # N132: was -18405.8

# Row 136
$ws.Cells.Item(136, 8).Value = 4731.439  # H136: was 4998.3687
$ws.Cells.Item(136, 9).Value = 2564.1614  # I136: was 2694.2144
$ws.Cells.Item(136, 11).Value = 7692.4842  # K136: was 8082.6432
$ws.Cells.Item(136, 13).Value = -5142.4842  # M136: was -5532.6432

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 1202.5652  # H3: was 1214.3684
$ws.Cells.Item(3, 9).Value = 1124.2222  # I3: was 1153
$ws.Cells.Item(3, 10).Value = 1484.6  # J3: was 1541.6666
$ws.Cells.Item(3, 11).Value = 1124.2222  # K3: was 1153
$ws.Cells.Item(3, 12).Value = 1484.6  # L3: was 1541.6666
$ws.Cells.Item(3, 13).Value = -1010.2222  # M3: was -1039
$ws.Cells.Item(3, 14).Value = -1712.6  # N3: was -1769.6666

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 2329.76  # H31: was 2334.2083
$ws.Cells.Item(31, 9).Value = 1445.4348  # I31: was 1410.091
$ws.Cells.Item(31, 11).Value = 1445.4348  # K31: was 1410.091
$ws.Cells.Item(31, 13).Value = -1150.4348  # M31: was -1115.091

# Row 34
$ws.Cells.Item(34, 8).Value = 2329.76  # H34: was 2334.2083
$ws.Cells.Item(34, 9).Value = 1445.4348  # I34: was 1410.091
$ws.Cells.Item(34, 11).Value = 1445.4348  # K34: was 1410.091
$ws.Cells.Item(34, 13).Value = -1243.4348  # M34: was -1208.091

# Row 51
$ws.Cells.Item(51, 8).Value = 53970  # H51: was 54166.332
$ws.Cells.Item(51, 10).Value = 59962.5  # J51: was 58999.6
$ws.Cells.Item(51, 12).Value = 59962.5  # L51: was 58999.6
$ws.Cells.Item(51, 14).Value = -61434.5  # N51: was -60471.6

# Row 61
$ws.Cells.Item(61, 8).Value = 53970  # H61: was 54166.332
$ws.Cells.Item(61, 10).Value = 59962.5  # J61: was 58999.6
$ws.Cells.Item(61, 12).Value = 59962.5  # L61: was 58999.6
$ws.Cells.Item(61, 14).Value = -60658.5  # N61: was -59695.6

# Row 99
$ws.Cells.Item(99, 8).Value = 4057.1428  # H99: was 4400
$ws.Cells.Item(99, 9).Value = 4050  # I99: was 4400
$ws.Cells.Item(99, 10).Value = 4100  # J99: was 0
$ws.Cells.Item(99, 11).Value = 4050  # K99: was 4400
$ws.Cells.Item(99, 12).Value = 4100  # L99: was 0
$ws.Cells.Item(99, 13).Value = -2552  # M99: was -2902
$ws.Cells.Item(99, 14).Value = -7096  # N99: was None

# Row 122
$ws.Cells.Item(122, 8).Value = 5936.5293  # H122: was 6150.8125
$ws.Cells.Item(122, 9).Value = 6175.273  # I122: was 6292.8
$ws.Cells.Item(122, 10).Value = 5498.8335  # J122: was 5914.1665
$ws.Cells.Item(122, 11).Value = 18525.819  # K122: was 18878.4
$ws.Cells.Item(122, 12).Value = 16496.5005  # L122: was 17742.4995
$ws.Cells.Item(122, 13).Value = -16075.819  # M122: was -16428.4
$ws.Cells.Item(122, 14).Value = -21396.5005  # N122: was -22642.4995

# Row 126
$ws.Cells.Item(126, 8).Value = 4057.1428  # H126: was 4400
$ws.Cells.Item(126, 9).Value = 4050  # I126: was 4400
$ws.Cells.Item(126, 10).Value = 4100  # J126: was 0
$ws.Cells.Item(126, 11).Value = 12150  # K126: was 13200
$ws.Cells.Item(126, 12).Value = 12300  # L126: was 0
$ws.Cells.Item(126, 13).Value = -9680  # M126: was -10730
$ws.Cells.Item(126, 14).Value = -17240  # N126: was None

# Row 132
$ws.Cells.Item(132, 8).Value = 1821331.6  # H132: was 2503536
$ws.Cells.Item(132, 9).Value = 2670073.8  # I132: was 4003977.5
$ws.Cells.Item(132, 10).Value = 2598.5715  # J132: was 2800
$ws.Cells.Item(132, 11).Value = 8010221.399999999  # K132: was 12011932.5
$ws.Cells.Item(132, 12).Value = 7795.7145  # L132: was 8400
$ws.Cells.Item(132, 13).Value = -8007691.399999999  # M132: was -12009402.5
$ws.Cells.Item(132, 14).Value = -12855.7145  # N132: was -13460

$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Cells.Item(80, 8).Value = 15000  # H80: was 0
$ws.Cells.Item(80, 10).Value = 15000  # J80: was 0
$ws.Cells.Item(80, 12).Value = 45000  # L80: was 0
$ws.Cells.Item(80, 14).Value = -46872  # N80: was None

# Row 83
$ws.Cells.Item(83, 8).Value = 15000  # H83: was 0
$ws.Cells.Item(83, 10).Value = 15000  # J83: was 0
$ws.Cells.Item(83, 12).Value = 135000  # L83: was 0
$ws.Cells.Item(83, 14).Value = -144360  # N83: was None

# Row 131
$ws.Cells.Item(131, 8).Value = 1574.5294  # H131: was 1615.381
$ws.Cells.Item(131, 10).Value = 1800.9166  # J131: was 1797.9375
$ws.Cells.Item(131, 12).Value = 5402.7498  # L131: was 5393.8125
$ws.Cells.Item(131, 14).Value = -15482.7498  # N131: was -15473.8125

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Cells.Item(2, 8).Value = 192.27272  # H2: was 224.86667
$ws.Cells.Item(2, 9).Value = 182.25  # I2: was 183.5
$ws.Cells.Item(2, 10).Value = 219  # J2: was 272.14285
$ws.Cells.Item(2, 11).Value = 182.25  # K2: was 183.5
$ws.Cells.Item(2, 12).Value = 219  # L2: was 272.14285
$ws.Cells.Item(2, 13).Value = -69.25  # M2: was -70.5
$ws.Cells.Item(2, 14).Value = -445  # N2: was -498.14285

# Row 12
$ws.Cells.Item(12, 8).Value = 30000  # H12: was 7001
$ws.Cells.Item(12, 10).Value = 30000  # J12: was 7001
$ws.Cells.Item(12, 12).Value = 30000  # L12: was 7001
$ws.Cells.Item(12, 14).Value = -30280  # N12: was -7281

# Row 41
$ws.Cells.Item(41, 8).Value = 929.6667  # H41: was 1499.5
$ws.Cells.Item(41, 9).Value = 929.6667  # I41: was 999.3333
$ws.Cells.Item(41, 10).Value = 0  # J41: was 3000
$ws.Cells.Item(41, 11).Value = 929.6667  # K41: was 999.3333
$ws.Cells.Item(41, 12).Value = 0  # L41: was 3000
$ws.Cells.Item(41, 13).Value = -574.6667  # M41: was -644.3333
$ws.Cells.Item(41, 14).ClearContents()  # N41: was -3710

# Row 122
$ws.Cells.Item(122, 8).Value = 1994.3334  # H122: was 1933.8889
$ws.Cells.Item(122, 9).Value = 1755.4762  # I122: was 1724.3043
$ws.Cells.Item(122, 10).Value = 3666.3333  # J122: was 2304.6924
$ws.Cells.Item(122, 11).Value = 5266.4286  # K122: was 5172.9129
$ws.Cells.Item(122, 12).Value = 10998.9999  # L122: was 6914.0772
$ws.Cells.Item(122, 13).Value = -2816.4286  # M122: was -2722.9129
$ws.Cells.Item(122, 14).Value = -15898.9999  # N122: was -11814.0772

# Row 126
$ws.Cells.Item(126, 8).Value = 3803.1  # H126: was 4048
$ws.Cells.Item(126, 9).Value = 3731.889  # I126: was 3998.5
$ws.Cells.Item(126, 11).Value = 11195.667  # K126: was 11995.5
$ws.Cells.Item(126, 13).Value = -8725.667000000001  # M126: was -9525.5

# Row 132
$ws.Cells.Item(132, 8).Value = 47622616  # H132: was 17547162
$ws.Cells.Item(132, 9).Value = 47622616  # I132: was 17547162
$ws.Cells.Item(132, 11).Value = 142867848  # K132: was 52641486
$ws.Cells.Item(132, 13).Value = -142865318  # M132: was -52638956

$ws = $wb.Worksheets.Item("LTW")
# Row 3
$ws.Cells.Item(3, 8).Value = 10005  # H3: was 1004
$ws.Cells.Item(3, 9).Value = 0  # I3: was 1004
$ws.Cells.Item(3, 10).Value = 10005  # J3: was 0
$ws.Cells.Item(3, 11).Value = 0  # K3: was 1004
$ws.Cells.Item(3, 12).Value = 10005  # L3: was 0
$ws.Cells.Item(3, 13).ClearContents()  # M3: was -892
$ws.Cells.Item(3, 14).Value = -10229  # N3: was None

# Row 15
$ws.Cells.Item(15, 8).Value = 10005  # H15: was 1004
$ws.Cells.Item(15, 9).Value = 0  # I15: was 1004
$ws.Cells.Item(15, 10).Value = 10005  # J15: was 0
$ws.Cells.Item(15, 11).Value = 0  # K15: was 1004
$ws.Cells.Item(15, 12).Value = 10005  # L15: was 0
$ws.Cells.Item(15, 13).ClearContents()  # M15: was -834
$ws.Cells.Item(15, 14).Value = -10345  # N15: was None

# Row 16
$ws.Cells.Item(16, 8).Value = 2001  # H16: was 2072.2856
$ws.Cells.Item(16, 9).Value = 1548.5  # I16: was 1648.3
$ws.Cells.Item(16, 11).Value = 1548.5  # K16: was 1648.3
$ws.Cells.Item(16, 13).Value = -1378.5  # M16: was -1478.3

# Row 61
$ws.Cells.Item(61, 8).Value = 2258.5  # H61: was 5930.923
$ws.Cells.Item(61, 9).Value = 2235.2  # I61: was 6577.4546
$ws.Cells.Item(61, 11).Value = 2235.2  # K61: was 6577.4546
$ws.Cells.Item(61, 13).Value = -2033.2  # M61: was -6375.4546

# Row 113
$ws.Cells.Item(113, 8).Value = 2258.5  # H113: was 5930.923
$ws.Cells.Item(113, 9).Value = 2235.2  # I113: was 6577.4546
$ws.Cells.Item(113, 11).Value = 2235.2  # K113: was 6577.4546
$ws.Cells.Item(113, 13).Value = -65.19999999999982  # M113: was -4407.4546

# Row 132
$ws.Cells.Item(132, 8).Value = 2581.5557  # H132: was 2645.758
$ws.Cells.Item(132, 9).Value = 2484.2642  # I132: was 2561.34
$ws.Cells.Item(132, 10).Value = 3097.2  # J132: was 2997.5
$ws.Cells.Item(132, 11).Value = 7452.792600000001  # K132: was 7684.02
$ws.Cells.Item(132, 12).Value = 9291.599999999999  # L132: was 8992.5
$ws.Cells.Item(132, 13).Value = -4922.792600000001  # M132: was -5154.02
$ws.Cells.Item(132, 14).Value = -14351.6  # N132: was -14052.5

# Row 136
$ws.Cells.Item(136, 8).Value = 4141.8335  # H136: was 4169.7393
$ws.Cells.Item(136, 9).Value = 3929.95  # I136: was 3952.5789
$ws.Cells.Item(136, 11).Value = 11789.85  # K136: was 11857.7367
$ws.Cells.Item(136, 13).Value = -9239.849999999999  # M136: was -9307.736699999999

$ws = $wb.Worksheets.Item("WVR")
# Row 7
$ws.Cells.Item(7, 8).Value = 6900  # H7: was 0
$ws.Cells.Item(7, 9).Value = 4350  # I7: was 0
$ws.Cells.Item(7, 10).Value = 12000  # J7: was 0
$ws.Cells.Item(7, 11).Value = 4350  # K7: was 0
$ws.Cells.Item(7, 12).Value = 12000  # L7: was 0
$ws.Cells.Item(7, 13).Value = -4237  # M7: was None
$ws.Cells.Item(7, 14).Value = -12226  # N7: was None

# Row 62
$ws.Cells.Item(62, 8).Value = 103159.125  # H62: was 135212.67
$ws.Cells.Item(62, 9).Value = 201194.25  # I62: was 161555.2
$ws.Cells.Item(62, 10).Value = 5124  # J62: was 3500
$ws.Cells.Item(62, 11).Value = 201194.25  # K62: was 161555.2
$ws.Cells.Item(62, 12).Value = 5124  # L62: was 3500
$ws.Cells.Item(62, 13).Value = -200570.25  # M62: was -160931.2
$ws.Cells.Item(62, 14).Value = -6372  # N62: was -4748

# Row 65
$ws.Cells.Item(65, 8).Value = 103159.125  # H65: was 135212.67
$ws.Cells.Item(65, 9).Value = 201194.25  # I65: was 161555.2
$ws.Cells.Item(65, 10).Value = 5124  # J65: was 3500
$ws.Cells.Item(65, 11).Value = 1005971.25  # K65: was 807776
$ws.Cells.Item(65, 12).Value = 25620  # L65: was 17500
$ws.Cells.Item(65, 13).Value = -1002851.25  # M65: was -804656
$ws.Cells.Item(65, 14).Value = -31860  # N65: was -23740

# Row 100
$ws.Cells.Item(100, 8).Value = 2002.5217  # H100: was 2019.9131
$ws.Cells.Item(100, 9).Value = 1607.5714  # I100: was 1626.619
$ws.Cells.Item(100, 11).Value = 3215.1428  # K100: was 3253.238
$ws.Cells.Item(100, 13).Value = -2674.1428  # M100: was -2712.238

# Row 132
$ws.Cells.Item(132, 8).Value = 2092.361  # H132: was 2178.2144
$ws.Cells.Item(132, 9).Value = 2092.361  # I132: was 2156
$ws.Cells.Item(132, 10).Value = 0  # J132: was 2311.5
$ws.Cells.Item(132, 11).Value = 6277.083  # K132: was 6468
$ws.Cells.Item(132, 12).Value = 0  # L132: was 6934.5
$ws.Cells.Item(132, 13).Value = -3747.083  # M132: was -3938
$ws.Cells.Item(132, 14).ClearContents()  # N132: was -11994.5
